$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Save the original values of the columns that change (D, J, K, L, M, P)
# for rows 2, 3, 4 before overwriting anything, since the update performs
# a cyclic rotation: row2<-row4, row3<-row2, row4<-row3.
$cols = @("D", "J", "K", "L", "M", "P")

$orig = @{}
foreach ($row in 2..4) {
    $orig[$row] = @{}
    foreach ($col in $cols) {
        $orig[$row][$col] = $ws.Range("$col$row").Value2
    }
}

# New row 2 = old row 4
foreach ($col in $cols) {
    $ws.Range("$col" + "2").Value = $orig[4][$col]
}

# New row 3 = old row 2
foreach ($col in $cols) {
    $ws.Range("$col" + "3").Value = $orig[2][$col]
}

# New row 4 = old row 3
foreach ($col in $cols) {
    $ws.Range("$col" + "4").Value = $orig[3][$col]
}
